$wb = $excel.ActiveWorkbook

# --- Insert a new "UseCaseContainer" sheet right after the "UseCase" sheet ---
$useCaseSheet = $wb.Worksheets.Item("UseCase")
$container1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $useCaseSheet)
$container1.Name = "UseCaseContainer"
$container1.Range("A1").Value = "container_name"
$container1.Range("B1").Value = "use_cases"

# --- Append a matching "UseCaseContainer1" sheet as the very last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$container2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$container2.Name = "UseCaseContainer1"
$container2.Range("A1").Value = "container_name"
$container2.Range("B1").Value = "use_cases"

# --- The newly appended last sheet becomes the active tab (activeTab moves from 29 to 31) ---
$container2.Activate()
